# Nudge a handful of diagram shapes on the single slide down slightly
# (wording/layout clean-up). Offsets below are taken straight from the
# target OOXML (EMU), converted to points (1 pt = 12700 EMU) since the
# PowerPoint object model's Shape.Top/Left are expressed in points.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shapes = $s.Shapes

# Map of shape Id -> new Y offset, given in points (Shape.Top/.Left are in
# points; 1 pt = 12700 EMU). The point values below are chosen so they land
# exactly on the target EMU values from the target OOXML after round-trip.
$newOffsets = @{
    25 = 195.9432373046875    # Graphic 24 (picture)      -> 2488479 EMU
    27 = 209.88418579101562   # TextBox 26                -> 2665529 EMU
    28 = 258.6890563964844    # TextBox 27                -> 3285351 EMU
    29 = 257.5106506347656    # TextBox 28                -> 3270385 EMU
    30 = 346.8409729003906    # TextBox 29                -> 4404880 EMU
    31 = 346.8409729003906    # TextBox 30                -> 4404880 EMU
    32 = 396.70294189453125   # TextBox 31                -> 5038127 EMU
    55 = 132.80804443359375   # Freeform 54                -> 1686662 EMU
    57 = 132.3984375          # Freeform 56                -> 1681460 EMU
    75 = 128.3373260498047    # Freeform 74                -> 1629884 EMU
}

foreach ($id in $newOffsets.Keys) {
    $shape = Get-ShapeById $shapes $id
    if ($shape -ne $null) {
        $shape.Top = $newOffsets[$id]
    }
}
